$d = $word.ActiveDocument

# The original text "sus honorario pactado..." needs to become
# "sus honorarios pactado..." (adding the missing "s" to "honorario").
$d.Content.Find.Execute("sus honorario pactado", $true, $false, $false, $false, $false,
                         $true, 1, $false, "sus honorarios pactado", 2)
